# Generate Report for Handoff
# Refresh handoff status/timestamps for file f56d6f71-23e5-4c7e-81dc-c62357c10949.md
# (now "Ready for handoff") and bump the "Latest Handoff" timestamps for the
# three files that were re-handed-off (rows for 4ca241b1..., f56d6f71..., fae5ebb3...)
# across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 6 -> 4ca241b1-d0db-4087-bb8a-2a1134fd6e31.md
$wsOverview.Range("D6").Value = "2016-20-20 06:20:48"

# Row 9 -> f56d6f71-23e5-4c7e-81dc-c62357c10949.md
$wsOverview.Range("B9").Value = "Ready for handoff"
$wsOverview.Range("C9").Value = "Ready for handoff"
$wsOverview.Range("D9").Value = "2016-20-20 06:20:48"

# Row 10 -> fae5ebb3-21dd-4a86-a96f-0a30b14d1210.md
$wsOverview.Range("D10").Value = "2016-20-20 06:20:48"

# ---- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 6 -> 4ca241b1-d0db-4087-bb8a-2a1134fd6e31.md
$wsZhCn.Range("E6").Value = "2016-03-20 06:20:44"

# Row 9 -> f56d6f71-23e5-4c7e-81dc-c62357c10949.md
$wsZhCn.Range("C9").Value = "Ready for handoff"
$wsZhCn.Range("E9").Value = "2016-03-20 06:20:44"

# Row 10 -> fae5ebb3-21dd-4a86-a96f-0a30b14d1210.md
$wsZhCn.Range("E10").Value = "2016-03-20 06:20:44"

# ---- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 6 -> 4ca241b1-d0db-4087-bb8a-2a1134fd6e31.md
$wsDeDe.Range("E6").Value = "2016-03-20 06:20:48"

# Row 9 -> f56d6f71-23e5-4c7e-81dc-c62357c10949.md
$wsDeDe.Range("C9").Value = "Ready for handoff"
$wsDeDe.Range("E9").Value = "2016-03-20 06:20:48"

# Row 10 -> fae5ebb3-21dd-4a86-a96f-0a30b14d1210.md
$wsDeDe.Range("E10").Value = "2016-03-20 06:20:48"
